$ws = $excel.ActiveWorkbook.ActiveSheet

# Cell updates derived from the cryptos list refresh.
# Columns B/C are plain text (coin name / link); D/E are forced to
# Text format before assignment so values like "0.9100" or "29.103.08"
# are preserved verbatim instead of being parsed as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.103.08'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.848.79'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.48%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7073'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -5.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.31'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9992'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3055'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07488'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.41'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08133'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7256'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.90%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.842.14'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.24'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.137.54'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.788'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -6.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '240.11'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007675'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.08'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9997'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.101.83'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.561'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1463'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -7.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.975'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.91'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.02'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.941'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.70%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -6.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.583'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.494'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.008'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05168'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.188'
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.034'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.26%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7076'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -7.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.641'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01864'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.675'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9100'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.992'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4297'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.067.18'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.42%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9991'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.33'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.756'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.072'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -7.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.182'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.37%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.75%  '
